# Commit: loss function changed to IoU + Center Distance
#
# Row 28 ("параметры теста 4" run): the note text loses its trailing
# period/space ("Размера батча изменён на 64. " -> "Размера батча изменён на 64").
# Two new log rows are appended describing the loss-function change:
#   row 29 - loss switched to plain IoU (with results + commit id)
#   row 30 - loss switched to a combination of IoU and center distance (newest entry, no results yet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the write order below intentionally mirrors the order the strings were
# first introduced in the authored workbook, so they land on the same
# shared-string indices (F29's text before F28's trimmed text, etc).

# --- Row 29: new entry - loss function replaced by IoU ---
$ws.Range("F29").Value = "Функция потерь заменена на IoU"

# --- Row 28: trim the trailing ". " from the batch-size note ---
$ws.Range("F28").Value = "Размера батча изменён на 64"

$ws.Range("H29").Value = "Train IoU: 0.43, Val IoU: 0.43, Test IoU: 0.40. Точность снизилась, но теперь модель более приближена к практической цели своей работы. "
$ws.Range("I29").Value = "3c89449"

# --- Row 30: new entry - loss function replaced by an IoU + center-distance combo ---
$ws.Range("F30").Value = "Функция потерь заменена на равносильную комбинацию IoU и расстояния между центрами"

$ws.Range("G29").Value = "параметры теста 4"
$ws.Range("G30").Value = "параметры теста 4"

$ws.Range("B29").Value = 1
$ws.Range("C29").Value = 40
$ws.Range("D29").Value = 19
$ws.Rows.Item(29).RowHeight = 60

$ws.Range("B30").Value = 1
$ws.Range("C30").Value = 40
$ws.Range("D30").Value = 20
$ws.Rows.Item(30).RowHeight = 45

# --- Move the active selection to the last-edited cell ---
$ws.Range("G30").Select() | Out-Null
